$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-shape the row layout.
#    Old rows: 4,6,7,9,11,12,13,15,16
#    New rows: 2,4,7,12,14,15,16,18,19  (row 12 is the blank "divider" row
#    that keeps moving down and keeps its special fill/format)
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).Delete()
$ws.Rows.Item(1).Delete()
$ws.Rows.Item(5).Resize(2).Insert()
$ws.Rows.Item(9).Resize(3).Insert()

# ---------------------------------------------------------------------------
# 2) New tip-percentage rows (5, 6, 8, 9, 10) next to the existing 15% (row 7)
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = "tip"
$ws.Range("D5").Value = 0.05
$ws.Range("D5").Style = "Percent"
$ws.Range("F5").Value = "CAPTURAR"

$ws.Range("C6").Value = "tip"
$ws.Range("D6").Value = 0.1
$ws.Range("D6").Style = "Percent"
$ws.Range("F6").Value = "CAPTURAR"

$ws.Range("C8").Value = "tip"
$ws.Range("D8").Value = 0.25
$ws.Range("D8").Style = "Percent"
$ws.Range("F8").Value = "CAPTURAR"

$ws.Range("C9").Value = "tip"
$ws.Range("D9").Value = 0.5
$ws.Range("D9").Style = "Percent"
$ws.Range("F9").Value = "CAPTURAR"

$ws.Range("C10").Value = "tip-custom"
$ws.Range("D10").Value = 0.3
$ws.Range("D10").Style = "Percent"
$ws.Range("F10").Value = "CAPTURAR"

# Remove the stray "todos os botões" note that used to sit next to the 15% row
$ws.Range("G7").ClearContents()

# ---------------------------------------------------------------------------
# 3) Relabel F18/F19 ("CALCULAR" -> "CALCULAR E MOSTRAR")
# ---------------------------------------------------------------------------
$ws.Range("F18").Value = "CALCULAR E MOSTRAR"
$ws.Range("F19").Value = "CALCULAR E MOSTRAR"

# ---------------------------------------------------------------------------
# 4) New rows 21/22 ("button reset" / "error message")
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = "button reset"
$ws.Range("F21").Value = "CAPTURAR"

$ws.Range("C22").Value = "error message"
$ws.Range("F22").Value = "CAPTURAR, innerhtml"

# ---------------------------------------------------------------------------
# 5) New column H, a running 1..12 index next to every "CAPTURAR"/"CALCULAR" row
# ---------------------------------------------------------------------------
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 2
$ws.Range("H6").Value = 3
$ws.Range("H7").Value = 4
$ws.Range("H8").Value = 5
$ws.Range("H9").Value = 6
$ws.Range("H10").Value = 7
$ws.Range("H12").Value = 8
$ws.Range("H14").Value = 9
$ws.Range("H16").Value = 10
$ws.Range("H21").Value = 11
$ws.Range("H22").Value = 12

# ---------------------------------------------------------------------------
# 6) Re-style every label cell in column C (except the "total tip"/"total
#    amount" rows, which stay as they were) with the new red, right-aligned
#    font used throughout the sheet.
# ---------------------------------------------------------------------------
$labelCells = "C4","C5","C6","C7","C8","C9","C10","C16","C18","C19","C21","C22"
foreach ($cellRef in $labelCells) {
    $cell = $ws.Range($cellRef)
    $cell.HorizontalAlignment = -4152
    $cell.Font.Bold = $false
    $cell.Font.Color = 255
}

# ---------------------------------------------------------------------------
# 7) Column G width + selection bookkeeping to match the saved view state
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 11.8
$ws.Range("C9").Select()
